$wb = $excel.ActiveWorkbook

# Both worksheets ("Table_1" and "Table_2") get a new row 1 inserted above the
# existing data. The old row 1 (the text column headers) becomes row 2, and
# every other row shifts down by one. The brand-new row 1 is a numeric
# index row (0,1,2,...,10) styled like the old header row.
$sheetNames = @("Table_1", "Table_2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Push every existing row down by one, leaving a blank row 1 behind.
    $ws.Rows.Item(1).Insert()

    # Copy the (now shifted-down) header row's formatting onto the new row 1
    # so it keeps the bold/centered/bordered header style.
    $ws.Range("A2:K2").Copy()
    $ws.Range("A1:K1").PasteSpecial(-4122)

    # Fill the new row 1 with the numeric column index values.
    $ws.Range("A1").Value = 0
    $ws.Range("B1").Value = 1
    $ws.Range("C1").Value = 2
    $ws.Range("D1").Value = 3
    $ws.Range("E1").Value = 4
    $ws.Range("F1").Value = 5
    $ws.Range("G1").Value = 6
    $ws.Range("H1").Value = 7
    $ws.Range("I1").Value = 8
    $ws.Range("J1").Value = 9
    $ws.Range("K1").Value = 10
}
